$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the startVal for row 2 (column B) from 0.5 to 0.3
$ws.Range("B2").Value = 0.3

# Update the selected cell to F5
$ws.Range("F5").Select()
